$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.295.80'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '2.616.28'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.13%  '
$c = $ws.Range("D5")
$c.Value = "'592.85"
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.52%  '
$c = $ws.Range("D6")
$c.Value = "'151.72"
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -0.10%  '
$c = $ws.Range("D9")
$c.Value = "'0.114"
$c.Style = "Normal"
$ws.Range('E9').Value = '  +4.86%  '
$c = $ws.Range("D10")
$c.Value = "'0.396"
$c.Style = "Normal"
$ws.Range('E10').Value = '  +3.57%  '
$c = $ws.Range("D11")
$c.Value = "'5.79"
$c.Style = "Normal"
$ws.Range('E11').Value = '  +1.75%  '
$c = $ws.Range("D12")
$c.Value = "'0.152"
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.19%  '
$c = $ws.Range("D13")
$c.Value = "'28.47"
$c.Style = "Normal"
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('D14').Value = '3.086.19'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '64.215.16'
$ws.Range('E15').Value = '  +1.24%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D16")
$c.Value = "'0.0000171"
$c.Style = "Normal"
$ws.Range('E16').Value = '  +12.18%  '
$ws.Range('D17').Value = '2.651.26'
$ws.Range('E17').Value = '  +1.75%  '
$c = $ws.Range("D18")
$c.Value = "'12.22"
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.77%  '
$c = $ws.Range("D19")
$c.Value = "'4.78"
$c.Style = "Normal"
$ws.Range('E19').Value = '  +2.07%  '
$c = $ws.Range("D20")
$c.Value = "'349.49"
$c.Style = "Normal"
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('E21').Value = '  +4.25%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('E23').Value = '  +1.32%  '
$ws.Range('E24').Value = '  -1.20%  '
$c = $ws.Range("D25")
$c.Value = "'9.27"
$c.Style = "Normal"
$ws.Range('E25').Value = '  +0.64%  '
$c = $ws.Range("D26")
$c.Value = "'1.64"
$c.Style = "Normal"
$ws.Range('E26').Value = '  -1.71%  '
$c = $ws.Range("D27")
$c.Value = "'8.23"
$c.Style = "Normal"
$ws.Range('E27').Value = '  +1.33%  '
$c = $ws.Range("D28")
$c.Value = "'0.163"
$c.Style = "Normal"
$ws.Range('E28').Value = '  +1.51%  '
$c = $ws.Range("D29")
$c.Value = "'545.07"
$c.Style = "Normal"
$ws.Range('E29').Value = '  -1.29%  '
$c = $ws.Range("D30")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').Value = '0.0₃0907'
$ws.Range('E31').Value = '  +7.25%  '
$c = $ws.Range("D32")
$c.Value = "'2.07"
$c.Style = "Normal"
$ws.Range('E32').Value = '  +1.43%  '
$c = $ws.Range("D33")
$c.Value = "'1.81"
$c.Style = "Normal"
$ws.Range('E33').Value = '  +3.02%  '
$c = $ws.Range("D34")
$c.Value = "'5.65"
$c.Style = "Normal"
$ws.Range('E34').Value = '  +8.08%  '
$c = $ws.Range("D35")
$c.Value = "'6.21"
$c.Style = "Normal"
$ws.Range('E35').Value = '  +0.81%  '
$c = $ws.Range("D36")
$c.Value = "'0.422"
$c.Style = "Normal"
$ws.Range('E36').Value = '  +2.11%  '
$c = $ws.Range("D37")
$c.Value = "'163.61"
$c.Style = "Normal"
$ws.Range('E37').Value = '  -2.43%  '
$c = $ws.Range("D38")
$c.Value = "'20.09"
$c.Style = "Normal"
$ws.Range('E38').Value = '  +3.15%  '
$c = $ws.Range("D39")
$c.Value = "'1.99"
$c.Style = "Normal"
$ws.Range('E39').Value = '  +3.21%  '
$c = $ws.Range("D40")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range('E40').Value = '  +0.01%  '
$c = $ws.Range("D42")
$c.Value = "'168.01"
$c.Style = "Normal"
$ws.Range('E42').Value = '  +0.81%  '
$c = $ws.Range("D43")
$c.Value = "'41.56"
$c.Style = "Normal"
$ws.Range('E43').Value = '  +4.67%  '
$c = $ws.Range("D44")
$c.Value = "'4.08"
$c.Style = "Normal"
$ws.Range('E44').Value = '  +4.58%  '
$c = $ws.Range("D45")
$c.Value = "'23.17"
$c.Style = "Normal"
$ws.Range('E45').Value = '  +7.11%  '
$c = $ws.Range("D46")
$c.Value = "'0.0596"
$c.Style = "Normal"
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('E47').Value = '  +10.68%  '
$c = $ws.Range("D48")
$c.Value = "'0.638"
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.67%  '
$c = $ws.Range("D49")
$c.Value = "'0.0250"
$c.Style = "Normal"
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  +1.53%  '
$c = $ws.Range("D51")
$c.Value = "'19.23"
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.26%  '
